$d = $word.ActiveDocument

# Update the date paragraph (first paragraph of the document)
$d.Paragraphs.Item(1).Range.Text = "2024-10-19 Saturday"

# Update the 100 math expressions in the table, cell by cell (row-major order)
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "99-77="
$t.Cell(1,2).Range.Text = "60-36="
$t.Cell(1,3).Range.Text = "54-50="
$t.Cell(1,4).Range.Text = "43-7="
$t.Cell(1,5).Range.Text = "65+6="
$t.Cell(2,1).Range.Text = "98-18="
$t.Cell(2,2).Range.Text = "20+31="
$t.Cell(2,3).Range.Text = "81-10="
$t.Cell(2,4).Range.Text = "91-8="
$t.Cell(2,5).Range.Text = "8+68="
$t.Cell(3,1).Range.Text = "48+17="
$t.Cell(3,2).Range.Text = "86-36="
$t.Cell(3,3).Range.Text = "61-61="
$t.Cell(3,4).Range.Text = "30-21="
$t.Cell(3,5).Range.Text = "73-37="
$t.Cell(4,1).Range.Text = "29+63="
$t.Cell(4,2).Range.Text = "39+16="
$t.Cell(4,3).Range.Text = "66-14="
$t.Cell(4,4).Range.Text = "4+94="
$t.Cell(4,5).Range.Text = "32+57="
$t.Cell(5,1).Range.Text = "17+10="
$t.Cell(5,2).Range.Text = "13+48="
$t.Cell(5,3).Range.Text = "69-62="
$t.Cell(5,4).Range.Text = "74-62="
$t.Cell(5,5).Range.Text = "16+77="
$t.Cell(6,1).Range.Text = "11+11="
$t.Cell(6,2).Range.Text = "64+23="
$t.Cell(6,3).Range.Text = "44-16="
$t.Cell(6,4).Range.Text = "75-28="
$t.Cell(6,5).Range.Text = "76+2="
$t.Cell(7,1).Range.Text = "81-23="
$t.Cell(7,2).Range.Text = "71-44="
$t.Cell(7,3).Range.Text = "33+59="
$t.Cell(7,4).Range.Text = "79-26="
$t.Cell(7,5).Range.Text = "16+67="
$t.Cell(8,1).Range.Text = "43+21="
$t.Cell(8,2).Range.Text = "18+31="
$t.Cell(8,3).Range.Text = "5+3="
$t.Cell(8,4).Range.Text = "13+2="
$t.Cell(8,5).Range.Text = "50-0="
$t.Cell(9,1).Range.Text = "66-63="
$t.Cell(9,2).Range.Text = "12+49="
$t.Cell(9,3).Range.Text = "48+50="
$t.Cell(9,4).Range.Text = "70-3="
$t.Cell(9,5).Range.Text = "94-53="
$t.Cell(10,1).Range.Text = "45+6="
$t.Cell(10,2).Range.Text = "90-4="
$t.Cell(10,3).Range.Text = "44-26="
$t.Cell(10,4).Range.Text = "93-59="
$t.Cell(10,5).Range.Text = "26-13="
$t.Cell(11,1).Range.Text = "84-54="
$t.Cell(11,2).Range.Text = "97-10="
$t.Cell(11,3).Range.Text = "21+17="
$t.Cell(11,4).Range.Text = "56-55="
$t.Cell(11,5).Range.Text = "62-17="
$t.Cell(12,1).Range.Text = "50-43="
$t.Cell(12,2).Range.Text = "46-28="
$t.Cell(12,3).Range.Text = "1+69="
$t.Cell(12,4).Range.Text = "76-28="
$t.Cell(12,5).Range.Text = "16+58="
$t.Cell(13,1).Range.Text = "84-37="
$t.Cell(13,2).Range.Text = "10+14="
$t.Cell(13,3).Range.Text = "91-69="
$t.Cell(13,4).Range.Text = "64-23="
$t.Cell(13,5).Range.Text = "16+80="
$t.Cell(14,1).Range.Text = "94-25="
$t.Cell(14,2).Range.Text = "31+21="
$t.Cell(14,3).Range.Text = "17-6="
$t.Cell(14,4).Range.Text = "0+68="
$t.Cell(14,5).Range.Text = "85-28="
$t.Cell(15,1).Range.Text = "8+8="
$t.Cell(15,2).Range.Text = "23+8="
$t.Cell(15,3).Range.Text = "81+5="
$t.Cell(15,4).Range.Text = "64-54="
$t.Cell(15,5).Range.Text = "75+24="
$t.Cell(16,1).Range.Text = "63-22="
$t.Cell(16,2).Range.Text = "21+50="
$t.Cell(16,3).Range.Text = "90-38="
$t.Cell(16,4).Range.Text = "35-3="
$t.Cell(16,5).Range.Text = "9+31="
$t.Cell(17,1).Range.Text = "86-50="
$t.Cell(17,2).Range.Text = "41+27="
$t.Cell(17,3).Range.Text = "16+70="
$t.Cell(17,4).Range.Text = "29-20="
$t.Cell(17,5).Range.Text = "83-35="
$t.Cell(18,1).Range.Text = "47-16="
$t.Cell(18,2).Range.Text = "26+8="
$t.Cell(18,3).Range.Text = "18+24="
$t.Cell(18,4).Range.Text = "36+1="
$t.Cell(18,5).Range.Text = "22-20="
$t.Cell(19,1).Range.Text = "27+12="
$t.Cell(19,2).Range.Text = "72-51="
$t.Cell(19,3).Range.Text = "23+44="
$t.Cell(19,4).Range.Text = "47+49="
$t.Cell(19,5).Range.Text = "6+71="
$t.Cell(20,1).Range.Text = "65-43="
$t.Cell(20,2).Range.Text = "75+13="
$t.Cell(20,3).Range.Text = "8+57="
$t.Cell(20,4).Range.Text = "43-34="
$t.Cell(20,5).Range.Text = "74+0="
